# Correct typos & update offloading fig
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the description cell (C2): clarify wording and punctuation.
$ws.Range("C2").Value = "Dictionary that contains input information such as {input node name: input tensor}."

# Update the active selection to C2 (matches the saved cursor position in the file).
$ws.Range("C2").Select()
